$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 1, shifting all existing rows down by one.
$ws.Rows.Item(1).Insert()

# The freshly inserted row 1 has no real content; copy the format/contents of
# the sheet's other blank filler row (now row 135) onto it so it is written
# out the same way as the sheet's other intentionally-blank rows.
$ws.Range("A135:D135").Copy($ws.Range("A1:D1"))

# Restore the sort definition (range grew by one row after the insert).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A114"))
$ws.Sort.SetRange($ws.Range("A2:D114"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Update the view: scroll position and selected cell.
$excel.ActiveWindow.ScrollRow = 93
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("G13").Select()
